$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the top of the data (just below the
# header row and the first, untouched, data row). This pushes the
# existing rows 4-51 down to rows 5-52, and the sheet grows from 51 to
# 52 data rows total (dimension A1:R52).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new observation. Columns
# that are constant across every other data row in this sheet are
# carried over unchanged.
$ws.Range("A4").Value = 8
$ws.Range("B4").Value = "Terminal La Palmera de La Serena"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 44756
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 100114007
$ws.Range("G4").Value = "Jengibre"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 14500
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 14750
$ws.Range("N4").Value = "`$/caja 13 kilos"
$ws.Range("O4").Value = "Perú"
$ws.Range("P4").Value = 1135
$ws.Range("Q4").Value = 13
$ws.Range("R4").Value = "Hortaliza"
